{"js": "const pairs = [\n  [\"630\u00f75=126, 0\", \"948\u00f73=316, 0\"],\n  [\"988\u00f74=247, 0\", \"779\u00f73=259, 2\"],\n  [\"299\u00f79=33, 2\", \"293\u00f74=73, 1\"],\n  [\"114\u00f74=28, 2\", \"328\u00f78=41, 0\"],\n  [\"677\u00f76=112, 5\", \"856\u00f76=142, 4\"],\n  [\"268\u00f76=44, 4\", \"825\u00f78=103, 1\"],\n  [\"573\u00f78=71, 5\", \"617\u00f79=68, 5\"],\n  [\"931\u00f77=133, 0\", \"600\u00f74=150, 0\"],\n  [\"943\u00f77=134, 5\", \"799\u00f74=199, 3\"],\n  [\"399\u00f79=44, 3\", \"856\u00f76=142, 4\"],\n  [\"108\u00f75=21, 3\", \"252\u00f76=42, 0\"],\n  [\"131\u00f74=32, 3\", \"657\u00f72=328, 1\"],\n  [\"120\u00f76=20, 0\", \"104\u00f77=14, 6\"],\n  [\"195\u00f78=24, 3\", \"154\u00f75=30, 4\"],\n  [\"890\u00f79=98, 8\", \"701\u00f77=100, 1\"],\n  [\"693\u00f72=346, 1\", \"474\u00f78=59, 2\"],\n  [\"630\u00f74=157, 2\", \"145\u00f77=20, 5\"],\n  [\"113\u00f77=16, 1\", \"142\u00f79=15, 7\"],\n  [\"610\u00f73=203, 1\", \"897\u00f76=149, 3\"],\n  [\"696\u00f77=99, 3\", \"992\u00f77=141, 5\"],\n  [\"465\u00f73=155, 0\", \"249\u00f75=49, 4\"],\n  [\"456\u00f78=57, 0\", \"392\u00f78=49, 0\"],\n  [\"765\u00f73=255, 0\", \"752\u00f76=125, 2\"],\n  [\"827\u00f78=103, 3\", \"484\u00f73=161, 1\"],\n  [\"911\u00f77=130, 1\", \"942\u00f79=104, 6\"],\n];\n\nfor (const [before, after] of pairs) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + before);\n  }\n  results.items[0].insertText(after, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n    @(\"630\u00f75=126, 0\", \"948\u00f73=316, 0\"),\n    @(\"988\u00f74=247, 0\", \"779\u00f73=259, 2\"),\n    @(\"299\u00f79=33, 2\", \"293\u00f74=73, 1\"),\n    @(\"114\u00f74=28, 2\", \"328\u00f78=41, 0\"),\n    @(\"677\u00f76=112, 5\", \"856\u00f76=142, 4\"),\n    @(\"268\u00f76=44, 4\", \"825\u00f78=103, 1\"),\n    @(\"573\u00f78=71, 5\", \"617\u00f79=68, 5\"),\n    @(\"931\u00f77=133, 0\", \"600\u00f74=150, 0\"),\n    @(\"943\u00f77=134, 5\", \"799\u00f74=199, 3\"),\n    @(\"399\u00f79=44, 3\", \"856\u00f76=142, 4\"),\n    @(\"108\u00f75=21, 3\", \"252\u00f76=42, 0\"),\n    @(\"131\u00f74=32, 3\", \"657\u00f72=328, 1\"),\n    @(\"120\u00f76=20, 0\", \"104\u00f77=14, 6\"),\n    @(\"195\u00f78=24, 3\", \"154\u00f75=30, 4\"),\n    @(\"890\u00f79=98, 8\", \"701\u00f77=100, 1\"),\n    @(\"693\u00f72=346, 1\", \"474\u00f78=59, 2\"),\n    @(\"630\u00f74=157, 2\", \"145\u00f77=20, 5\"),\n    @(\"113\u00f77=16, 1\", \"142\u00f79=15, 7\"),\n    @(\"610\u00f73=203, 1\", \"897\u00f76=149, 3\"),\n    @(\"696\u00f77=99, 3\", \"992\u00f77=141, 5\"),\n    @(\"465\u00f73=155, 0\", \"249\u00f75=49, 4\"),\n    @(\"456\u00f78=57, 0\", \"392\u00f78=49, 0\"),\n    @(\"765\u00f73=255, 0\", \"752\u00f76=125, 2\"),\n    @(\"827\u00f78=103, 3\", \"484\u00f73=161, 1\"),\n    @(\"911\u00f77=130, 1\", \"942\u00f79=104, 6\"),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
